$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.298.01'
$ws.Range('E2').Value = '  +1.33%  '
$ws.Range('D3').Value = '2.251.01'
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.42'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.09'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.95%  '
$ws.Range('E7').Value = '  +1.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.530'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('E10').Value = '  -1.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0818'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.04%  '
$ws.Range('E12').Value = '  -0.71%  '
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('D14').Value = '2.594.21'
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').Value = '2.331.26'
$ws.Range('E15').Value = '  +4.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.838'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.67'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.31%  '
$ws.Range('D18').Value = '44.139.42'
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').Value = '0.0₃0974'
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.25'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -5.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.40'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.83%  '
$ws.Range('E22').Value = '  +1.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.33'
$ws.Range('D23').ClearFormats()
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('E25').Value = '  +0.39%  '
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.06'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.04%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.19'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.66%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.21'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.00'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.24'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '153.18'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0803'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.76%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.61'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.40%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.22'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.120'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +3.36%  '
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('E38').Value = '  -5.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.62'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '14.76'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.87%  '
$ws.Range('E41').Value = '  -3.24%  '
$ws.Range('E42').Value = '  -1.80%  '
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').Value = '1.756.15'
$ws.Range('E44').Value = '  +3.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '83.78'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '100.84'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.98'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.20'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '55.13'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.57'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.23%  '
